# Apply the changes described in the commit:
#  - 執行梯次 (column C) values change from "第N梯次" wording to "第N次"
#  - 狀態 (column G) values are updated to reflect new progress state
#    (進行中/待開始 -> 已完成/進行中/未開始)
#  - 狀態 column is widened to better fit the new, longer status badges
#  - Active selection / cursor position moves from K6 to G5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: 執行梯次 wording update ---
$ws.Range("C2").Value = "第一次"
$ws.Range("C3").Value = "第一次"
$ws.Range("C4").Value = "第二次"
$ws.Range("C5").Value = "第二次"
$ws.Range("C6").Value = "第三次"
$ws.Range("C7").Value = "第三次"

# --- Column G: 狀態 update ---
$ws.Range("G2").Value = "已完成"
$ws.Range("G3").Value = "已完成"
$ws.Range("G4").Value = "進行中"
$ws.Range("G5").Value = "未開始"
$ws.Range("G6").Value = "未開始"
$ws.Range("G7").Value = "未開始"

# --- Widen the 狀態 column so the new badge text fits ---
$ws.Columns("G").ColumnWidth = 30.49

# --- Update the saved cursor/selection position ---
$ws.Range("G5").Select()
